{"js": "// Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" paragraph,\n// the \"\u00a9 2020 . Contact: ... Creative Commons Attribution\" paragraph that\n// follows it, and the (now orphaned) blank paragraph that used to sit\n// between the copyright line and the trailing page-break paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the two text paragraphs by their content so this does not depend\n// on a hard-coded paragraph index.\nlet verNoIndex = -1;\nlet copyrightIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (verNoIndex === -1 && t.indexOf(\"Ver no Jupiter\") !== -1) {\n    verNoIndex = i;\n  }\n  if (copyrightIndex === -1 && t.indexOf(\"Powered by Jekyll and Github pages\") !== -1) {\n    copyrightIndex = i;\n  }\n}\n\nif (verNoIndex === -1 || copyrightIndex === -1) {\n  throw new Error(\"Could not locate the footer paragraphs to remove.\");\n}\n\n// The blank paragraph that immediately follows the copyright paragraph\n// (if any) gets removed along with it, since the edit collapses the two\n// trailing blank-paragraph separators down to one.\nlet blankAfterCopyrightIndex = -1;\nif (copyrightIndex + 1 < items.length && items[copyrightIndex + 1].text === \"\") {\n  blankAfterCopyrightIndex = copyrightIndex + 1;\n}\n\n// Delete from the highest index down to the lowest so earlier indices stay\n// valid while we work.\nconst toDelete = [blankAfterCopyrightIndex, copyrightIndex, verNoIndex]\n  .filter(i => i !== -1)\n  .sort((a, b) => b - a);\n\nfor (const idx of toDelete) {\n  items[idx].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" paragraph,\n# the \"\u00a9 2020 . Contact: ... Creative Commons Attribution\" paragraph that\n# follows it, and the (now orphaned) blank paragraph that used to sit\n# between the copyright line and the trailing page-break paragraph.\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$verNoIndex = -1\n$copyrightIndex = -1\n\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($verNoIndex -eq -1 -and $t -like \"*Ver no Jupiter*\") {\n        $verNoIndex = $i\n    }\n    if ($copyrightIndex -eq -1 -and $t -like \"*Powered by Jekyll and Github pages*\") {\n        $copyrightIndex = $i\n    }\n}\n\nif ($verNoIndex -eq -1 -or $copyrightIndex -eq -1) {\n    throw \"Could not locate the footer paragraphs to remove.\"\n}\n\n# The blank paragraph immediately after the copyright paragraph (if any) gets\n# removed along with it, since the edit collapses the two trailing\n# blank-paragraph separators down to one.\n$blankAfterCopyrightIndex = -1\nif (($copyrightIndex + 1) -le $count) {\n    $nextText = $d.Paragraphs.Item($copyrightIndex + 1).Range.Text\n    if ($nextText -eq \"`r\") {\n        $blankAfterCopyrightIndex = $copyrightIndex + 1\n    }\n}\n\n# Delete from the highest index down to the lowest so earlier indices stay\n# valid while we work.\n$indices = @($blankAfterCopyrightIndex, $copyrightIndex, $verNoIndex) |\n    Where-Object { $_ -ne -1 } |\n    Sort-Object -Descending -Unique\n\nforeach ($idx in $indices) {\n    $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
